# "Mostrar valores na tela 3" - add a 3rd sheet (Planilha3) with sample
# values, refresh the linear-system inputs / solution on Planilha2 and
# show them in (scientific) exponential notation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add a new worksheet "Planilha3" at the end of the workbook and
#    fill it with the sample matrix values.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Planilha3"

$ws3Data = @(
    @(3, 4, 8, 6),
    @(9, 3, 6, 4),
    @(5, 2, 1, 5)
)
for ($i = 0; $i -lt $ws3Data.Length; $i++) {
    for ($j = 0; $j -lt $ws3Data[$i].Length; $j++) {
        $ws3.Cells.Item($i + 2, $j + 1).Value = $ws3Data[$i][$j]
    }
}

# ---------------------------------------------------------------------
# 2) Planilha2: update the system's coefficients (A/C/E/G columns of
#    rows 19-21) and the solution vector (B25:D25) with new values, and
#    show the solution in scientific notation.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Planilha2")

$ws2.Range("B25").Value = [double]"9.9556898207726e56"
$ws2.Range("C25").Value = [double]"9.9556898207726e56"
$ws2.Range("D25").Value = [double]"-3.98227592830904e57"
$ws2.Range("B25:D25").NumberFormat = "0.00E+00"

$ws2.Range("A19").Value = 3
$ws2.Range("C19").Value = 2
$ws2.Range("E19").Value = -1
$ws2.Range("G19").Value = 0

$ws2.Range("A20").Value = 1
$ws2.Range("C20").Value = 3
$ws2.Range("E20").Value = 1
$ws2.Range("G20").Value = 1

$ws2.Range("A21").Value = 2
$ws2.Range("C21").Value = 2
$ws2.Range("E21").Value = -2
$ws2.Range("G21").Value = 2

# ---------------------------------------------------------------------
# 3) Restore Planilha2 as the active sheet/selection (A10) so the view
#    matches what the author left behind.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A10").Select()
